$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI ligand/receptor-expressing-cell counts (1 -> 3) for rows 2-10,
# together with every downstream expression/specificity metric recomputed from the
# refreshed counts (per "Natmi following Dr Hou advice").
# Each row entry: row number, then a map of column letter -> new value.
$updates = @(
    @{ Row = 2; Cells = @{ "E" = 3; "G" = 5.138644333333333; "H" = 15.415933; "I" = 0.1788080791399461; "J" = 0.1788080791399461; "K" = 3; "M" = 29.880108; "N" = 89.64032399999999; "O" = 0.4868991363731112; "P" = 0.4868991363731112; "Q" = 153.543247653588; "R" = 1381.889228882292; "S" = 0.08706149930977465; "T" = 0.08706149930977466 } },
    @{ Row = 3; Cells = @{ "E" = 3; "G" = 5.138644333333333; "H" = 15.415933; "I" = 0.1788080791399461; "J" = 0.1788080791399461; "K" = 3; "M" = 25.29401133333333; "N" = 75.882034; "O" = 0.412168265041468; "P" = 0.412168265041468; "Q" = 129.9769280053024; "R" = 1169.792352047722; "S" = 0.07369901575450907; "T" = 0.07369901575450907 } },
    @{ Row = 4; Cells = @{ "E" = 3; "G" = 5.138644333333333; "H" = 15.415933; "I" = 0.1788080791399461; "J" = 0.1788080791399461; "K" = 3; "M" = 6.194048666666666; "N" = 18.582146; "O" = 0.1009325985854208; "P" = 0.1009325985854208; "Q" = 31.82901308135755; "R" = 286.4611177322179; "S" = 0.01804756407566233; "T" = 0.01804756407566233 } },
    @{ Row = 5; Cells = @{ "E" = 3; "G" = 18.194752; "H" = 54.584256; "I" = 0.6331180841693511; "J" = 0.6331180841693511; "K" = 3; "M" = 29.880108; "N" = 89.64032399999999; "O" = 0.4868991363731112; "P" = 0.4868991363731112; "Q" = 543.6611547932159; "R" = 4892.950393138944; "S" = 0.3082646484042558; "T" = 0.3082646484042558 } },
    @{ Row = 6; Cells = @{ "E" = 3; "G" = 18.194752; "H" = 54.584256; "I" = 0.6331180841693511; "J" = 0.6331180841693511; "K" = 3; "M" = 25.29401133333333; "N" = 75.882034; "O" = 0.412168265041468; "P" = 0.412168265041468; "Q" = 460.2182632951894; "R" = 4141.964369656705; "S" = 0.2609511823184595; "T" = 0.2609511823184596 } },
    @{ Row = 7; Cells = @{ "E" = 3; "G" = 18.194752; "H" = 54.584256; "I" = 0.6331180841693511; "J" = 0.6331180841693511; "K" = 3; "M" = 6.194048666666666; "N" = 18.582146; "O" = 0.1009325985854208; "P" = 0.1009325985854208; "Q" = 112.6991793659307; "R" = 1014.292614293376; "S" = 0.06390225344663578; "T" = 0.06390225344663578 } },
    @{ Row = 8; Cells = @{ "E" = 3; "G" = 5.404926666666667; "H" = 16.21478; "I" = 0.1880738366907028; "J" = 0.1880738366907027; "K" = 3; "M" = 29.880108; "N" = 89.64032399999999; "O" = 0.4868991363731112; "P" = 0.4868991363731112; "Q" = 161.49979253208; "R" = 1453.49813278872; "S" = 0.09157298865908073; "T" = 0.09157298865908071 } },
    @{ Row = 9; Cells = @{ "E" = 3; "G" = 5.404926666666667; "H" = 16.21478; "I" = 0.1880738366907028; "J" = 0.1880738366907027; "K" = 3; "M" = 25.29401133333333; "N" = 75.882034; "O" = 0.412168265041468; "P" = 0.412168265041468; "Q" = 136.7122763625023; "R" = 1230.41048726252; "S" = 0.07751806696849935; "T" = 0.07751806696849935 } },
    @{ Row = 10; Cells = @{ "E" = 3; "G" = 5.404926666666667; "H" = 16.21478; "I" = 0.1880738366907028; "J" = 0.1880738366907027; "K" = 3; "M" = 6.194048666666666; "N" = 18.582146; "O" = 0.1009325985854208; "P" = 0.1009325985854208; "Q" = 33.47837881309778; "R" = 301.30540931788; "S" = 0.01898278106312269; "T" = 0.01898278106312269 } }
)

foreach ($entry in $updates) {
    $r = $entry.Row
    foreach ($col in $entry.Cells.Keys) {
        $ws.Range("$col$r").Value = $entry.Cells[$col]
    }
}